$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "62.609.53", "0.507")
# that Excel would silently coerce to a Double via a plain .Value assignment,
# losing the exact text (and turning some into scientific notation).
# Force the cell to Text format first, assign, then restore the default
# "Normal" style so no stray style index is left behind on the cell.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '62.609.53'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +2.75%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.943.14'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '591.30'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.57%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '148.88'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +7.51%  '
$ws.Range('E7').Value = '  -0.04%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.507'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.04%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.942.20'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.89%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.10'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.83%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.151'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +10.06%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.438'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +2.61%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000236'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +8.88%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '32.47'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  -0.37%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.430.85'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.94%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '62.617.22'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +2.98%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '6.66'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.93%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '2.946.43'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.78%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '438.06'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +3.24%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.48'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.26%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.665'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.99%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.98'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.25'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +8.99%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '80.38'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.72%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.91'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +4.82%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.12'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +3.81%  '
$ws.Range('E28').Value = '  +0.02%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.33'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +10.95%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0000103'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +23.49%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '2.59'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +2.66%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +5.33%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.110'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.42%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '26.14'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +2.26%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.10%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.988'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.11'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +11.35%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '5.58'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.04%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '49.65'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.43%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.02'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +6.93%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '8.41'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.116'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.87%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.278'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +5.30%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '40.32'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +6.78%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.705.63'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.66%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '135.51'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +3.35%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0341'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.98%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '356.50'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  +2.41%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '22.76'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +2.74%  '
